$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: hours went from 2 to 0.5, and a new comment/note is added in G11 ---
$ws.Range("E11").Value = 0.5

# New G11 cell: note about the finished task / image, using the same
# center/vcenter/wrap formatting used elsewhere for "description"-style cells
# (e.g. C9, D3, D9 ...). Copy the formatting from C9 (same style) instead of
# setting alignment properties one-by-one, to avoid generating throwaway
# intermediate cell styles.
$ws.Range("C9").Copy()
$g11 = $ws.Range("G11")
$g11.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$g11.Value = "fini … un peu random comme image mais fonctionnelle, le temps est GOOD"

# --- Extend the "Terminé" conditional formatting (currently covering
# F2, A3:F19, G7, G5) so it also covers the new G11 cell. ---
$cf = $g11.FormatConditions.Add(9, 0, "Terminé")
$cf.Text = "Terminé"
$cf.Formula1 = 'NOT(ISERROR(SEARCH("Terminé",A2)))'
$cf.Font.Color = 393372
$cf.Interior.Color = 13551615

# --- Move the active selection cursor from B12 to D14 ---
$ws.Range("D14").Select()
